$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 61149.715
$ws.Range("I19").Value = 899.5
$ws.Range("J19").Value = 85249.8
$ws.Range("K19").Value = 899.5
$ws.Range("L19").Value = 85249.8
$ws.Range("M19").Value = -724.5
$ws.Range("N19").Value = -85599.8

$ws.Range("H33").Value = 7813543
$ws.Range("I33").Value = 8929032
$ws.Range("J33").Value = 5121.75
$ws.Range("K33").Value = 8929032
$ws.Range("L33").Value = 5121.75
$ws.Range("M33").Value = -8928803
$ws.Range("N33").Value = -5579.75

$ws.Range("H40").Value = 2922.75
$ws.Range("I40").Value = 1348
$ws.Range("J40").Value = 4497.5
$ws.Range("K40").Value = 1348
$ws.Range("L40").Value = 4497.5
$ws.Range("M40").Value = -1173
$ws.Range("N40").Value = -4847.5

$ws.Range("H69").Value = 9620.52
$ws.Range("I69").Value = 7200
$ws.Range("J69").Value = 11854.846
$ws.Range("K69").Value = 21600
$ws.Range("L69").Value = 35564.538
$ws.Range("M69").Value = -20726
$ws.Range("N69").Value = -37312.538

$ws.Range("H72").Value = 9620.52
$ws.Range("I72").Value = 7200
$ws.Range("J72").Value = 11854.846
$ws.Range("K72").Value = 64800
$ws.Range("L72").Value = 106693.614
$ws.Range("M72").Value = -60432
$ws.Range("N72").Value = -115429.614

$ws.Range("H76").Value = 6161.5454
$ws.Range("I76").Value = 6042.5557
$ws.Range("J76").Value = 6697
$ws.Range("K76").Value = 6042.5557
$ws.Range("L76").Value = 6697
$ws.Range("M76").Value = -5727.5557
$ws.Range("N76").Value = -7327

$ws.Range("H79").Value = 6161.5454
$ws.Range("I79").Value = 6042.5557
$ws.Range("J79").Value = 6697
$ws.Range("K79").Value = 6042.5557
$ws.Range("L79").Value = 6697
$ws.Range("M79").Value = -4950.5557
$ws.Range("N79").Value = -8881

$ws.Range("H100").Value = 1793.9231
$ws.Range("I100").Value = 1368.3334
$ws.Range("J100").Value = 2751.5
$ws.Range("K100").Value = 1368.3334
$ws.Range("L100").Value = 2751.5
$ws.Range("M100").Value = -827.3334
$ws.Range("N100").Value = -3833.5

$ws.Range("H112").Value = 3821.9473
$ws.Range("I112").Value = 1266.6666
$ws.Range("J112").Value = 4301.0625
$ws.Range("K112").Value = 3799.9998
$ws.Range("L112").Value = 12903.1875
$ws.Range("M112").Value = -2691.9998
$ws.Range("N112").Value = -15119.1875

$ws.Range("H138").Value = 2387.61
$ws.Range("I138").Value = 1191.4762
$ws.Range("J138").Value = 2705.5696
$ws.Range("K138").Value = 3574.4286
$ws.Range("L138").Value = 8116.708799999999
$ws.Range("M138").Value = 1565.5714
$ws.Range("N138").Value = -18396.7088

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1215.619
$ws.Range("I97").Value = 1185.5
$ws.Range("J97").Value = 1312
$ws.Range("K97").Value = 1185.5
$ws.Range("L97").Value = 1312
$ws.Range("M97").Value = -689.5
$ws.Range("N97").Value = -2304

$ws.Range("H102").Value = 17242.312
$ws.Range("I102").Value = 1875.3334
$ws.Range("J102").Value = 36999.855
$ws.Range("K102").Value = 1875.3334
$ws.Range("L102").Value = 36999.855
$ws.Range("M102").Value = -253.3334
$ws.Range("N102").Value = -40243.855

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 10120.833
$ws.Range("I99").Value = 1225
$ws.Range("J99").Value = 19016.666
$ws.Range("K99").Value = 1225
$ws.Range("L99").Value = 19016.666
$ws.Range("M99").Value = 273
$ws.Range("N99").Value = -22012.666

$ws.Range("H105").Value = 1772.5555
$ws.Range("I105").Value = 790.8
$ws.Range("J105").Value = 2999.75
$ws.Range("K105").Value = 790.8
$ws.Range("L105").Value = 2999.75
$ws.Range("M105").Value = 956.2
$ws.Range("N105").Value = -6493.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 6651.7393
$ws.Range("I99").Value = 1999
$ws.Range("J99").Value = 7094.857
$ws.Range("K99").Value = 1999
$ws.Range("L99").Value = 7094.857
$ws.Range("M99").Value = -501
$ws.Range("N99").Value = -10090.857

$ws.Range("H126").Value = 6651.7393
$ws.Range("I126").Value = 1999
$ws.Range("J126").Value = 7094.857
$ws.Range("K126").Value = 5997
$ws.Range("L126").Value = 21284.571
$ws.Range("M126").Value = -3527
$ws.Range("N126").Value = -26224.571

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1462.27
$ws.Range("I131").Value = 789.6
$ws.Range("J131").Value = 1497.6737
$ws.Range("K131").Value = 2368.8
$ws.Range("L131").Value = 4493.0211
$ws.Range("M131").Value = 2671.2
$ws.Range("N131").Value = -14573.0211

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1019.73914
$ws.Range("I2").Value = 1394.3125
$ws.Range("J2").Value = 163.57143
$ws.Range("K2").Value = 1394.3125
$ws.Range("L2").Value = 163.57143
$ws.Range("M2").Value = -1281.3125
$ws.Range("N2").Value = -389.57143

$ws.Range("H9").Value = 691.5
$ws.Range("I9").Value = 1825
$ws.Range("J9").Value = 124.75
$ws.Range("K9").Value = 1825
$ws.Range("L9").Value = 124.75
$ws.Range("M9").Value = -1655
$ws.Range("N9").Value = -464.75

$ws.Range("H11").Value = 633849.2
$ws.Range("I11").Value = 1431286.4
$ws.Range("J11").Value = 13620.223
$ws.Range("K11").Value = 1431286.4
$ws.Range("L11").Value = 13620.223
$ws.Range("M11").Value = -1431147.4
$ws.Range("N11").Value = -13898.223

$ws.Range("H18").Value = 1504166.4
$ws.Range("I18").Value = 1801399.6
$ws.Range("J18").Value = 18000
$ws.Range("K18").Value = 1801399.6
$ws.Range("L18").Value = 18000
$ws.Range("M18").Value = -1801106.6
$ws.Range("N18").Value = -18586

$ws.Range("H19").Value = 1034.375
$ws.Range("I19").Value = 1000
$ws.Range("J19").Value = 1055
$ws.Range("K19").Value = 1000
$ws.Range("L19").Value = 1055
$ws.Range("M19").Value = -712
$ws.Range("N19").Value = -1631

$ws.Range("H20").Value = 15342.481
$ws.Range("I20").Value = 2000
$ws.Range("J20").Value = 26016.467
$ws.Range("K20").Value = 2000
$ws.Range("L20").Value = 26016.467
$ws.Range("M20").Value = -1755
$ws.Range("N20").Value = -26506.467

$ws.Range("H53").Value = 23555.555
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 23555.555
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 23555.555
$ws.Range("M53").ClearContents()
$ws.Range("N53").Value = -24817.555

$ws.Range("H63").Value = 74990
$ws.Range("J63").Value = 74990
$ws.Range("L63").Value = 74990
$ws.Range("N63").Value = -76362

$ws.Range("H66").Value = 74990
$ws.Range("J66").Value = 74990
$ws.Range("L66").Value = 224970
$ws.Range("N66").Value = -231834

$ws.Range("H80").Value = 16207.947
$ws.Range("I80").Value = 11165
$ws.Range("J80").Value = 27134.334
$ws.Range("K80").Value = 11165
$ws.Range("L80").Value = 27134.334
$ws.Range("M80").Value = -10167
$ws.Range("N80").Value = -29130.334

$ws.Range("H83").Value = 16207.947
$ws.Range("I83").Value = 11165
$ws.Range("J83").Value = 27134.334
$ws.Range("K83").Value = 55825
$ws.Range("L83").Value = 135671.67
$ws.Range("M83").Value = -50833
$ws.Range("N83").Value = -145655.67

$ws.Range("H87").Value = 24000
$ws.Range("J87").Value = 24000
$ws.Range("L87").Value = 24000
$ws.Range("N87").Value = -26496

$ws.Range("H88").Value = 150130
$ws.Range("J88").Value = 150130
$ws.Range("L88").Value = 150130
$ws.Range("N88").Value = -151032

$ws.Range("H90").Value = 24000
$ws.Range("J90").Value = 24000
$ws.Range("L90").Value = 72000
$ws.Range("N90").Value = -84480

$ws.Range("H91").Value = 150130
$ws.Range("J91").Value = 150130
$ws.Range("L91").Value = 150130
$ws.Range("N91").Value = -153250

$ws.Range("H94").Value = 21998
$ws.Range("J94").Value = 21998
$ws.Range("L94").Value = 21998
$ws.Range("N94").Value = -23350

$ws.Range("H95").Value = 64207.168
$ws.Range("J95").Value = 64207.168
$ws.Range("L95").Value = 64207.168
$ws.Range("N95").Value = -69699.16800000001

$ws.Range("H97").Value = 4220.1763
$ws.Range("I97").Value = 1871.8462
$ws.Range("J97").Value = 11852.25
$ws.Range("K97").Value = 1871.8462
$ws.Range("L97").Value = 11852.25
$ws.Range("M97").Value = -1375.8462
$ws.Range("N97").Value = -12844.25

$ws.Range("H99").Value = 9192.5
$ws.Range("I99").Value = 5090.1665
$ws.Range("J99").Value = 21499.5
$ws.Range("K99").Value = 5090.1665
$ws.Range("L99").Value = 21499.5
$ws.Range("M99").Value = -2844.1665
$ws.Range("N99").Value = -25991.5

$ws.Range("H101").Value = 65700
$ws.Range("J101").Value = 65700
$ws.Range("L101").Value = 65700
$ws.Range("N101").Value = -72190

$ws.Range("H102").Value = 463362.38
$ws.Range("I102").Value = 509258.25
$ws.Range("J102").Value = 4403.5
$ws.Range("K102").Value = 509258.25
$ws.Range("L102").Value = 4403.5
$ws.Range("M102").Value = -507636.25
$ws.Range("N102").Value = -7647.5

$ws.Range("H105").Value = 178316.5
$ws.Range("J105").Value = 178316.5
$ws.Range("L105").Value = 178316.5
$ws.Range("N105").Value = -185304.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3308.6875
$ws.Range("I22").Value = 2895.3333
$ws.Range("J22").Value = 4548.75
$ws.Range("K22").Value = 2895.3333
$ws.Range("L22").Value = 4548.75
$ws.Range("M22").Value = -2600.3333
$ws.Range("N22").Value = -5138.75

$ws.Range("H27").Value = 3308.6875
$ws.Range("I27").Value = 2895.3333
$ws.Range("J27").Value = 4548.75
$ws.Range("K27").Value = 2895.3333
$ws.Range("L27").Value = 4548.75
$ws.Range("M27").Value = -2788.3333
$ws.Range("N27").Value = -4762.75

$ws.Range("H46").Value = 2700.6956
$ws.Range("I46").Value = 1219.2858
$ws.Range("J46").Value = 3348.8125
$ws.Range("K46").Value = 1219.2858
$ws.Range("L46").Value = 3348.8125
$ws.Range("M46").Value = -1031.2858
$ws.Range("N46").Value = -3724.8125

$ws.Range("H93").Value = 16202.454
$ws.Range("I93").Value = 25699.75
$ws.Range("J93").Value = 10775.429
$ws.Range("K93").Value = 25699.75
$ws.Range("L93").Value = 10775.429
$ws.Range("M93").Value = -24451.75
$ws.Range("N93").Value = -13271.429

$ws.Range("H100").Value = 3482.1904
$ws.Range("I100").Value = 4856.4165
$ws.Range("J100").Value = 1649.8889
$ws.Range("K100").Value = 4856.4165
$ws.Range("L100").Value = 1649.8889
$ws.Range("M100").Value = -4315.4165
$ws.Range("N100").Value = -2731.8889

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4328
$ws.Range("I81").Value = 2672.5
$ws.Range("J81").Value = 5431.6665
$ws.Range("K81").Value = 5345
$ws.Range("L81").Value = 10863.333
$ws.Range("M81").Value = -4284
$ws.Range("N81").Value = -12985.333

$ws.Range("H84").Value = 4328
$ws.Range("I84").Value = 2672.5
$ws.Range("J84").Value = 5431.6665
$ws.Range("K84").Value = 26725
$ws.Range("L84").Value = 54316.665
$ws.Range("M84").Value = -21421
$ws.Range("N84").Value = -64924.665
